$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.974.57"
$ws.Range("E2").Value = "  -5.23%  "

$ws.Range("D3").Value = "3.369.85"
$ws.Range("E3").Value = "  -7.21%  "

$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.45%  "

$ws.Range("D5").Value = "'183.46"
$ws.Range("E5").Value = "  -8.89%  "

$ws.Range("D6").Value = "'527.96"
$ws.Range("E6").Value = "  -9.60%  "

$ws.Range("D7").Value = "'0.597"
$ws.Range("E7").Value = "  -4.37%  "

$ws.Range("D8").Value = "3.362.98"
$ws.Range("E8").Value = "  -7.25%  "

$ws.Range("E9").Value = "  -0.04%  "

$ws.Range("D10").Value = "'0.617"
$ws.Range("E10").Value = "  -10.42%  "

$ws.Range("D11").Value = "'56.94"
$ws.Range("E11").Value = "  -6.32%  "

$ws.Range("D12").Value = "'0.131"
$ws.Range("E12").Value = "  -13.15%  "

$ws.Range("D13").Value = "'0.0000251"
$ws.Range("E13").Value = "  -12.43%  "

$ws.Range("D14").Value = "'9.14"
$ws.Range("E14").Value = "  -9.84%  "

$ws.Range("D15").Value = "3.886.48"
$ws.Range("E15").Value = "  -7.67%  "

$ws.Range("E16").Value = "  -4.18%  "

$ws.Range("D17").Value = "3.360.25"
$ws.Range("E17").Value = "  -7.67%  "

$ws.Range("D18").Value = "64.521.05"
$ws.Range("E18").Value = "  -5.75%  "

$ws.Range("D19").Value = "'17.26"
$ws.Range("E19").Value = "  -10.35%  "

$ws.Range("D20").Value = "'10.95"
$ws.Range("E20").Value = "  -13.06%  "

$ws.Range("D21").Value = "'0.958"
$ws.Range("E21").Value = "  -11.27%  "

$ws.Range("D22").Value = "'370.85"
$ws.Range("E22").Value = "  -8.89%  "

$ws.Range("D23").Value = "'80.67"
$ws.Range("E23").Value = "  -6.16%  "

$ws.Range("D24").Value = "'3.69"
$ws.Range("E24").Value = "  -14.46%  "

$ws.Range("D25").Value = "'10.68"
$ws.Range("E25").Value = "  -18.18%  "

$ws.Range("D26").Value = "'3.73"
$ws.Range("E26").Value = "  -6.88%  "

$ws.Range("D27").Value = "'5.87"
$ws.Range("E27").Value = "  -4.66%  "

$ws.Range("D28").Value = "'2.62"
$ws.Range("E28").Value = "  -11.00%  "

$ws.Range("D29").Value = "'11.23"
$ws.Range("E29").Value = "  -12.03%  "

$ws.Range("D30").Value = "'8.35"
$ws.Range("E30").Value = "  -11.99%  "

$ws.Range("D31").Value = "'666.08"
$ws.Range("E31").Value = "  -2.44%  "

$ws.Range("D32").Value = "'29.03"
$ws.Range("E32").Value = "  -8.89%  "

$ws.Range("D33").Value = "'6.70"
$ws.Range("E33").Value = "  -14.65%  "

$ws.Range("B34").Value = "Cosmos"
$ws.Range("C34").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D34").Value = "'11.04"
$ws.Range("E34").Value = "  -10.43%  "

$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").Value = "'60.30"
$ws.Range("E35").Value = "  -5.69%  "

$ws.Range("E36").Value = "  -10.08%  "

$ws.Range("D38").Value = "'36.06"
$ws.Range("E38").Value = "  -14.40%  "

$ws.Range("E39").Value = "  -11.13%  "

$ws.Range("D40").Value = "'0.995"
$ws.Range("E40").Value = "  -0.40%  "

$ws.Range("E41").Value = "  -7.10%  "

$ws.Range("D42").Value = "2.794.74"
$ws.Range("E42").Value = "  -12.56%  "

$ws.Range("D43").Value = "'2.73"
$ws.Range("E43").Value = "  -15.72%  "

$ws.Range("D44").Value = "'2.58"
$ws.Range("E44").Value = "  -10.18%  "

$ws.Range("D45").Value = "0.0₃0611"
$ws.Range("E45").Value = "  -21.34%  "

$ws.Range("D46").Value = "'0.0385"
$ws.Range("E46").Value = "  -8.59%  "

$ws.Range("D47").Value = "'2.30"
$ws.Range("E47").Value = "  -15.23%  "

$ws.Range("D48").Value = "'0.124"
$ws.Range("E48").Value = "  -6.06%  "

$ws.Range("D49").Value = "'135.58"
$ws.Range("E49").Value = "  -2.07%  "

$ws.Range("D50").Value = "'2.81"
$ws.Range("E50").Value = "  -9.17%  "

$ws.Range("D51").Value = "'2.57"
$ws.Range("E51").Value = "  -6.49%  "
